$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings, added in this exact order so they land at the ---
# --- same shared-string indices (24, 25, 26) as the target workbook.   ---
$ws.Range("A12").Value = "Engineering (`$/Wdc): Changed to be ""Developer Overhead (`$/Wdc)"""
$ws.Range("A15").Value = "Land preparation: Changed to be ""Transmission Line (`$/Wdc)"""
$ws.Range("A5").Value = "Inverter (`$/Wdc)"

# --- Row 4 ---
$ws.Range("C4").Value = 0.69666666666666699
$ws.Range("E4").Value = 0.68
$ws.Range("G4").Value = 0.65

# --- Row 5 ---
$ws.Range("C5").Value = 0.28999999999999998
$ws.Range("E5").Value = 0.13
$ws.Range("G5").Value = 0.11

# --- Row 6 ---
$ws.Range("C6").Formula = "=0.329105205267494+0.203586501897641"
$ws.Range("E6").Formula = "=0.178+0.156+0.029"
$ws.Range("G6").Formula = "=0.16+0.16"

# --- Row 7 ---
$ws.Range("C7").Formula = "=0.215199265236845+0.115106555743999"
$ws.Range("E7").Formula = "=0.55*0.187+0.45*0.187"
$ws.Range("G7").Value = 0.19

# --- Row 8 ---
$ws.Range("C8").Formula = "=0.677401030654478+0.37787140321983"
$ws.Range("E8").Formula = "=0.06+0.1661/2+0.1661/2+0.493"
$ws.Range("G8").Formula = "=0.103724570767137+0.06"

# --- Row 9 ---
$ws.Range("C9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("G9").Value = 0

# --- Row 11 ---
$ws.Range("C11").Value = 0.11
$ws.Range("E11").Value = 0.001
$ws.Range("G11").Value = 0.023767500513328502

# --- Row 12 ---
$ws.Range("G12").Formula = "=0.16"

# --- Row 13 ---
$ws.Range("F13").Value = 0.03
$ws.Range("G13").Value = 0.03

# --- Row 14 ---
$ws.Range("G14").Value = 0.03

# --- Row 15 ---
$ws.Range("G15").Value = 0.02

# --- Row 16 ---
$ws.Range("C16").Value = 0.51679868381138649
$ws.Range("E16").Value = 0.67240632966186697
$ws.Range("G16").Value = 1

# --- Row 17 ---
$ws.Range("C17").Value = 0.05
$ws.Range("E17").Value = 0.05
$ws.Range("G17").Value = 0.05

# --- Row 18 (totals) ---
$ws.Range("C18").Formula = "=SUM(C4:C8)*(1+C9)*(1+C16*C17)+SUM(C11:C15)"
$ws.Range("E18").Formula = "=SUM(E4:E8)*(1+E9)*(1+E16*E17)+SUM(E11:E15)"
$ws.Range("G18").Formula = "=SUM(G4:G8)*(1+G9)*(1+G16*G17)+SUM(G11:G15)"

# --- New currency-with-fill number format (style 25) for the "existing"/ ---
# --- "proposed" columns F & G that previously only carried the fill.    ---
$fg = $ws.Range("F4:G8,F11:G15")
$fg.NumberFormat = """$""#,##0.00"

# --- Highlight the two relabeled line items with a yellow fill (style 26) ---
$ws.Range("A12").Interior.Color = 65535
$ws.Range("A15").Interior.Color = 65535

# --- Give C18 its own currency format distinct from the shared B18 one ---
# --- (style 27).                                                        ---
$ws.Range("C18").NumberFormat = "_(""$""* #,##0.00_);_(""$""* (#,##0.00);_(""$""* ""-""??_);_(@_)x"
$ws.Range("C18").NumberFormat = "_(""$""* #,##0.00_);_(""$""* \(#,##0.00\);_(""$""* ""-""??_);_(@_)"

# --- Mark rows as having an explicit (custom) height, matching their ---
# --- existing height so nothing visually changes.                   ---
foreach ($r in 4,5,6,7,8,9,10,11,16) {
    $ws.Rows.Item($r).RowHeight = 14.45
}

# --- Move the active selection, as captured by the saved view state. ---
$ws.Range("A5").Select()
